$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 32 (pushes the old Total/notes rows down to 33-35)
$ws.Rows("32:32").Insert()

# Fill in the new expense line: Home Depot Order for Wire Strippers
$ws.Range("A32").Value = "Home Depot Order"
$ws.Range("B32").Value = 43191
$ws.Range("C32").Value = "Morgan"
$ws.Range("D32").Value = "Home Depot Order 1.pdf"
$ws.Range("E32").Value = 30.43
$ws.Range("F32").Value = "Wire Strippers"

# Give D32 the Hyperlink style + add the actual hyperlink (like the other receipt cells)
$ws.Hyperlinks.Add($ws.Range("D32"), "Home%20Depot%20Order%201.pdf", "", "", "Home Depot Order 1.pdf")
$ws.Range("D32").Style = "Hyperlink"

# Update the Total Expenses formula (now in row 33) to include the new row
$ws.Range("E33").Formula = "=SUM(E2:E32)"

# Restore the view/selection state from the authored edit
$ws.Range("E33").Select()
